# Delete the duplicate "PENDAHULUAN" outline slide (slide 2 of 8).
# Everything else in the deck is left untouched; PowerPoint repacks the
# slide id list / relationship ids (and drops the now-unused "Open Sans
# Bold" embedded font, since it was only referenced on this slide) as a
# natural consequence of removing the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.Delete()
